# --- Update the "Control" sheet's data table with new columns/rows -------
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Control")

# Resize the existing table (Table26, A1:J8) out to M8 so we have 13
# columns (3 new blank ones appended at K/L/M for now).
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:M8"))

# --- Re-point column headers into their new target layout ----------------
# Final layout:
#  A index | B source | C citation_doi | D var_code | E variable_name |
#  F description | G year_min | H year_max | I frequency | J granularity |
#  K country_coverage | L drive_link | M notes
$ws.Range("A1").Value2 = "index"
$ws.Range("B1").Value2 = "source"
$ws.Range("C1").Value2 = "citation_doi"
$ws.Range("D1").Value2 = "var_code"
$ws.Range("E1").Value2 = "variable_name"
$ws.Range("F1").Value2 = "description"
$ws.Range("G1").Value2 = "year_min"
$ws.Range("H1").Value2 = "year_max"
$ws.Range("I1").Value2 = "frequency"
$ws.Range("J1").Value2 = "granularity"
$ws.Range("K1").Value2 = "country_coverage"
$ws.Range("L1").Value2 = "drive_link"
$ws.Range("M1").Value2 = "notes"

# --- Data rows -------------------------------------------------------------
$ws.Range("A2").Value2 = 1
$ws.Range("B2").Value2 = "Bank Danych Lokalnych (GUS"
$ws.Range("C2").Value2 = "link to source"
$ws.Range("F2").Value2 = "Gross domestic product per capita"
$ws.Range("G2").Value2 = 2000
$ws.Range("H2").Value2 = 2023
$ws.Range("I2").Value2 = "annual"
$ws.Range("J2").Value2 = "NUTS-3"
$ws.Range("K2").Value2 = "Poland"
$ws.Range("L2").Value2 = "link to folder"

$ws.Range("A3").Value2 = 2
$ws.Range("B3").Value2 = "Bank Danych Lokalnych (GUS"
$ws.Range("C3").Value2 = "link to source"
$ws.Range("F3").Value2 = "Gross domestic product per capita, Poland=100"
$ws.Range("G3").Value2 = 2000
$ws.Range("H3").Value2 = 2023
$ws.Range("I3").Value2 = "annual"
$ws.Range("J3").Value2 = "NUTS-4"
$ws.Range("K3").Value2 = "Poland"
$ws.Range("L3").Value2 = "link to folder"

$ws.Range("A4").Value2 = 3
$ws.Range("B4").Value2 = "Bank Danych Lokalnych (GUS"
$ws.Range("C4").Value2 = "link to source"
$ws.Range("F4").Value2 = "Gross domestic product per capita, region=100"
$ws.Range("G4").Value2 = 2000
$ws.Range("H4").Value2 = 2023
$ws.Range("I4").Value2 = "annual"
$ws.Range("J4").Value2 = "NUTS-5"
$ws.Range("K4").Value2 = "Poland"
$ws.Range("L4").Value2 = "link to folder"

$ws.Range("A5").Value2 = 4
$ws.Range("B5").Value2 = "EUROSTAT"
$ws.Range("C5").Value2 = "link to source"
$ws.Range("F5").Value2 = "Geofiles"
$ws.Range("H5").Value2 = 2024
$ws.Range("J5").Value2 = "NUTS-3"
$ws.Range("K5").Value2 = "EU"
$ws.Range("L5").Value2 = "link to folder"

# --- Hyperlinks (order matches the authored workbook: L2, L3, L4, C2,
#     C3:C4, L5, C5 so the relationship ids line up rId1..rId7, with the
#     table part becoming rId8) -------------------------------------------
$ws.Hyperlinks.Add($ws.Range("L2"), "https://drive.google.com/drive/folders/gus-gdp-per-capita") | Out-Null
$ws.Hyperlinks.Add($ws.Range("L3"), "https://drive.google.com/drive/folders/gus-gdp-per-capita") | Out-Null
$ws.Hyperlinks.Add($ws.Range("L4"), "https://drive.google.com/drive/folders/gus-gdp-per-capita") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://bdl.stat.gov.pl/bdl/start") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3:C4"), "https://bdl.stat.gov.pl/bdl/start", "", "", "link to source") | Out-Null
$ws.Hyperlinks.Add($ws.Range("L5"), "https://drive.google.com/drive/folders/eurostat-geofiles") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"), "https://ec.europa.eu/eurostat/web/main/data/database") | Out-Null

# make sure all of the hyperlink cells pick up the "Hyperlink" cell style
# (Hyperlinks.Add only auto-styles the anchor cell of each call)
$ws.Range("C2").Style = "Hyperlink"
$ws.Range("C3:C4").Style = "Hyperlink"
$ws.Range("C5").Style = "Hyperlink"
$ws.Range("L2").Style = "Hyperlink"
$ws.Range("L3").Style = "Hyperlink"
$ws.Range("L4").Style = "Hyperlink"
$ws.Range("L5").Style = "Hyperlink"

# --- Wrap text formatting on the new "description" column ------------------
# (order matters for which new cell-style index each variant gets: the
# wrap-on style must be minted before the wrap-off/plain-alignment style)
$ws.Range("F5:F8").WrapText = $true
$ws.Range("F1:F4").WrapText = $false

# --- View tweaks -------------------------------------------------------
$excel.ActiveWindow.Zoom = 110
$ws.Range("C5").Select() | Out-Null
